$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New batch of generated "Proctor" test records that replace the contents
# of the data rows (rows 2-31) on the "Bulk Upload" sheet.
$data = @(
    @("Proctor97998", "Automation97998", "proctorautomation97998@gmail.com", "97998"),
    @("Proctor76074", "Automation76074", "proctorautomation76074@gmail.com", "76074"),
    @("Proctor53252", "Automation53252", "proctorautomation53252@gmail.com", "53252"),
    @("Proctor78247", "Automation78247", "proctorautomation78247@gmail.com", "78247"),
    @("Proctor64363", "Automation64363", "proctorautomation64363@gmail.com", "64363"),
    @("Proctor59592", "Automation59592", "proctorautomation59592@gmail.com", "59592"),
    @("Proctor40110", "Automation40110", "proctorautomation40110@gmail.com", "40110"),
    @("Proctor78443", "Automation78443", "proctorautomation78443@gmail.com", "78443"),
    @("Proctor80491", "Automation80491", "proctorautomation80491@gmail.com", "80491"),
    @("Proctor46067", "Automation46067", "proctorautomation46067@gmail.com", "46067"),
    @("Proctor77233", "Automation77233", "proctorautomation77233@gmail.com", "77233"),
    @("Proctor77490", "Automation77490", "proctorautomation77490@gmail.com", "77490"),
    @("Proctor48199", "Automation48199", "proctorautomation48199@gmail.com", "48199"),
    @("Proctor15436", "Automation15436", "proctorautomation15436@gmail.com", "15436"),
    @("Proctor86281", "Automation86281", "proctorautomation86281@gmail.com", "86281"),
    @("Proctor54836", "Automation54836", "proctorautomation54836@gmail.com", "54836"),
    @("Proctor00471", "Automation00471", "proctorautomation00471@gmail.com", "00471"),
    @("Proctor66885", "Automation66885", "proctorautomation66885@gmail.com", "66885"),
    @("Proctor24591", "Automation24591", "proctorautomation24591@gmail.com", "24591"),
    @("Proctor16625", "Automation16625", "proctorautomation16625@gmail.com", "16625"),
    @("Proctor99913", "Automation99913", "proctorautomation99913@gmail.com", "99913"),
    @("Proctor10454", "Automation10454", "proctorautomation10454@gmail.com", "10454"),
    @("Proctor00511", "Automation00511", "proctorautomation00511@gmail.com", "00511"),
    @("Proctor90070", "Automation90070", "proctorautomation90070@gmail.com", "90070"),
    @("Proctor87319", "Automation87319", "proctorautomation87319@gmail.com", "87319"),
    @("Proctor63714", "Automation63714", "proctorautomation63714@gmail.com", "63714"),
    @("Proctor92552", "Automation92552", "proctorautomation92552@gmail.com", "92552"),
    @("Proctor77271", "Automation77271", "proctorautomation77271@gmail.com", "77271"),
    @("Proctor72658", "Automation72658", "proctorautomation72658@gmail.com", "72658"),
    @("Proctor78061", "Automation78061", "proctorautomation78061@gmail.com", "78061")
)

# Column D holds numeric-looking ID strings (some zero-padded, e.g. "00471").
# Force text storage so leading zeros survive, without leaving a lasting
# number-format/style change on the cells.
$idRange = $ws.Range("D2:D31")
$idRange.NumberFormat = "@"

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}

$idRange.ClearFormats()
